# Update "want to go" counts (column F) on several sheets, each value
# incremented by 1 except two which decrease by 1 (488 -> 487).
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1748
$ws1.Range("F5").Value = 439
$ws1.Range("F12").Value = 1385
$ws1.Range("F14").Value = 317
$ws1.Range("F15").Value = 653
$ws1.Range("F16").Value = 12605
$ws1.Range("F17").Value = 12637
$ws1.Range("F22").Value = 44
$ws1.Range("F23").Value = 487
$ws1.Range("F24").Value = 1976

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 11
$ws2.Range("F10").Value = 66

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 150

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 150
$ws4.Range("F6").Value = 1748
$ws4.Range("F7").Value = 439
$ws4.Range("F17").Value = 1385
$ws4.Range("F19").Value = 317
$ws4.Range("F21").Value = 653
$ws4.Range("F22").Value = 12605
$ws4.Range("F23").Value = 12637
$ws4.Range("F28").Value = 44
$ws4.Range("F29").Value = 487
$ws4.Range("F30").Value = 11
$ws4.Range("F32").Value = 1976
$ws4.Range("F39").Value = 66
